$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("I1:I1048576")
$rng.FormatConditions.Delete()
$fc1 = $rng.FormatConditions.Add(9, 0, [Type]::Missing, [Type]::Missing, "FAILED", 0)
Write-Host "count:" $rng.FormatConditions.Count
Write-Host "type:" $fc1.Type
Write-Host "text:" $fc1.Text
